$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at 992-994 (shifts existing rows 992:1032 down to 995:1035,
# growing the used range from A1:T1032 to A1:T1035) to hold a new weekly
# price-report block for Comercializadora del Agro de Limari - Limon.
$ws.Rows("992:994").Insert()

# Row 992: 1a amarillo
$ws.Range("A992").Value = 2
$ws.Range("B992").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C992").Value = "Coquimbo"
$ws.Range("D992").Value = 45267
$ws.Range("E992").Value = 4
$ws.Range("F992").Value = "Fruta"
$ws.Range("G992").Value = 100102
$ws.Range("H992").Value = "Cítricos"
$ws.Range("I992").Value = 100102003
$ws.Range("J992").Value = "Limón"
$ws.Range("K992").Value = "Sin especificar"
$ws.Range("L992").Value = "1a amarillo"
$ws.Range("M992").Value = 900
$ws.Range("N992").Value = 9300
$ws.Range("O992").Value = 9500
$ws.Range("P992").Value = 9400
$ws.Range("Q992").Value = "$/malla 18 kilos"
$ws.Range("R992").Value = "Provincia de Limarí"
$ws.Range("S992").Value = 522
$ws.Range("T992").Value = 18

# Row 993: 2a amarillo
$ws.Range("A993").Value = 2
$ws.Range("B993").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C993").Value = "Coquimbo"
$ws.Range("D993").Value = 45267
$ws.Range("E993").Value = 4
$ws.Range("F993").Value = "Fruta"
$ws.Range("G993").Value = 100102
$ws.Range("H993").Value = "Cítricos"
$ws.Range("I993").Value = 100102003
$ws.Range("J993").Value = "Limón"
$ws.Range("K993").Value = "Sin especificar"
$ws.Range("L993").Value = "2a amarillo"
$ws.Range("M993").Value = 750
$ws.Range("N993").Value = 6300
$ws.Range("O993").Value = 6500
$ws.Range("P993").Value = 6400
$ws.Range("Q993").Value = "$/malla 18 kilos"
$ws.Range("R993").Value = "Provincia de Limarí"
$ws.Range("S993").Value = 356
$ws.Range("T993").Value = 18

# Row 994: 3a amarillo
$ws.Range("A994").Value = 2
$ws.Range("B994").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C994").Value = "Coquimbo"
$ws.Range("D994").Value = 45267
$ws.Range("E994").Value = 4
$ws.Range("F994").Value = "Fruta"
$ws.Range("G994").Value = 100102
$ws.Range("H994").Value = "Cítricos"
$ws.Range("I994").Value = 100102003
$ws.Range("J994").Value = "Limón"
$ws.Range("K994").Value = "Sin especificar"
$ws.Range("L994").Value = "3a amarillo"
$ws.Range("M994").Value = 600
$ws.Range("N994").Value = 4300
$ws.Range("O994").Value = 4500
$ws.Range("P994").Value = 4400
$ws.Range("Q994").Value = "$/malla 18 kilos"
$ws.Range("R994").Value = "Provincia de Limarí"
$ws.Range("S994").Value = 244
$ws.Range("T994").Value = 18
